$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsZh.Range("E2").Value = "2016-03-19 14:54:07"
$wsZh.Range("H2").Value = "2016-03-19 14:54:25"

# de-de sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsDe.Range("E2").Value = "2016-03-19 14:54:10"
$wsDe.Range("H2").Value = "2016-03-19 14:54:31"
